$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# New row 7 - text/label cells
$ws.Range("A7").Value = "CW3M"
$ws.Range("B7").Value = "Baseline 2010-18 C374"
$ws.Range("C7").Value = "2010-18"

# New row 7 - numeric cells
$ws.Range("D7").Value = 577.95190099999991
$ws.Range("E7").Value = 2094.2995878888887
$ws.Range("F7").Value = 5.8562380000000012
$ws.Range("G7").Value = 232.21855144444442
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 6.5062423333333328
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 559.18725922222222
$ws.Range("L7").Value = 89.974237444444441
$ws.Range("M7").Value = 1651.4325493333336
$ws.Range("N7").Value = 616.21212433333346
$ws.Range("O7").Value = 15529.390950555557
$ws.Range("P7").Value = 2216.8192002222222
$ws.Range("Q7").Value = -0.026350555555555506
$ws.Range("R7").Value = -0.00010888888888888886

# Number formats matching the rest of the table
$ws.Range("D7").NumberFormat = "0.00"
$ws.Range("E7:G7").NumberFormat = "0.00"
$ws.Range("H7:J7").NumberFormat = "0.00"
$ws.Range("K7:L7").NumberFormat = "0.00"
$ws.Range("M7").NumberFormat = "0.00"
$ws.Range("N7").NumberFormat = "0.00"
$ws.Range("O7:P7").NumberFormat = "0"
$ws.Range("Q7").NumberFormat = "0.00"
$ws.Range("R7").NumberFormat = "0.000000"

# Highlighted (yellow fill) cells matching the other changed-value cells in the table
$ws.Range("D7").Interior.Color = 65535
$ws.Range("K7").Interior.Color = 65535
$ws.Range("L7").Interior.Color = 65535
$ws.Range("N7").Interior.Color = 65535

# Update the active selection to the newly added cell
$null = $ws.Range("D7").Select()
